$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "Execute"
$ws.Range("A9").Value = "Ejecución "

$ws.Range("A10").Value = "Enlace"
$ws.Range("B10").Value = "Link"

$ws.Range("B11").Value = "Design"
$ws.Range("A11").Value = "Diseño"

$ws.Range("A12").Value = "Group "
$ws.Range("B12").Value = "grupo "

$ws.Range("A13").Select()
